# The commit swaps the two theme parts shipped with this deck:
#   ppt/theme/theme1.xml (bound to the slide master, i.e. the deck's
#   visible design) goes from the colourful "Integral" / "Red Violet"
#   theme to the plain default "Office Theme" / "Office" colour scheme.
#
# PowerPoint's object model exposes a theme's 12 colour scheme slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink, in that fixed order)
# as Master.Theme.ThemeColorScheme.Item(1..12).RGB, so we repoint the
# slide master's theme colours at the stock "Office" palette.
#
# Note: ThemeColorScheme.Item(n).RGB takes a VBA-style RGB() long, i.e.
# the bytes are packed as 0xBBGGRR (red in the low byte), not 0xRRGGBB.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0x000000   # dk1      -> 000000
$cs.Item(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$cs.Item(3).RGB  = 0x6A5444   # dk2      -> 44546A
$cs.Item(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$cs.Item(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$cs.Item(6).RGB  = 0x317DED   # accent2  -> ED7D31
$cs.Item(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$cs.Item(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$cs.Item(9).RGB  = 0xC47244   # accent5  -> 4472C4
$cs.Item(10).RGB = 0x47AD70   # accent6  -> 70AD47
$cs.Item(11).RGB = 0xC16305   # hlink    -> 0563C1
$cs.Item(12).RGB = 0x724F95   # folHlink -> 954F72
